$d = $word.ActiveDocument

# Update the date/title paragraph (first paragraph, outside the table)
$d.Paragraphs.Item(1).Range.Text = "2025-12-24 Wednesday"

# Update each table cell value, addressed positionally (row, col) to
# avoid ambiguity since some old values repeat (e.g. "5+39=" appears twice).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50-27="
$t.Cell(1, 2).Range.Text = "92-16="
$t.Cell(1, 3).Range.Text = "24+69="
$t.Cell(1, 4).Range.Text = "36+38="
$t.Cell(1, 5).Range.Text = "82-13="

$t.Cell(2, 1).Range.Text = "53+29="
$t.Cell(2, 2).Range.Text = "60-38="
$t.Cell(2, 3).Range.Text = "42-6="
$t.Cell(2, 4).Range.Text = "21-15="
$t.Cell(2, 5).Range.Text = "56-37="

$t.Cell(3, 1).Range.Text = "57+4="
$t.Cell(3, 2).Range.Text = "94-35="
$t.Cell(3, 3).Range.Text = "10-2="
$t.Cell(3, 4).Range.Text = "44+8="
$t.Cell(3, 5).Range.Text = "24-5="

$t.Cell(4, 1).Range.Text = "62-19="
$t.Cell(4, 2).Range.Text = "14+19="
$t.Cell(4, 3).Range.Text = "72-65="
$t.Cell(4, 4).Range.Text = "39+14="
$t.Cell(4, 5).Range.Text = "42-24="

$t.Cell(5, 1).Range.Text = "41-7="
$t.Cell(5, 2).Range.Text = "96-27="
$t.Cell(5, 3).Range.Text = "39+8="
$t.Cell(5, 4).Range.Text = "81-18="
$t.Cell(5, 5).Range.Text = "69+16="

$t.Cell(6, 1).Range.Text = "3+49="
$t.Cell(6, 2).Range.Text = "58-49="
$t.Cell(6, 3).Range.Text = "24+38="
$t.Cell(6, 4).Range.Text = "6+5="
$t.Cell(6, 5).Range.Text = "83-15="

$t.Cell(7, 1).Range.Text = "97-89="
$t.Cell(7, 2).Range.Text = "70-5="
$t.Cell(7, 3).Range.Text = "73-5="
$t.Cell(7, 4).Range.Text = "56+25="
$t.Cell(7, 5).Range.Text = "76-17="

$t.Cell(8, 1).Range.Text = "31-28="
$t.Cell(8, 2).Range.Text = "63+8="
$t.Cell(8, 3).Range.Text = "3+18="
$t.Cell(8, 4).Range.Text = "52-18="
$t.Cell(8, 5).Range.Text = "29+59="

$t.Cell(9, 1).Range.Text = "83-69="
$t.Cell(9, 2).Range.Text = "56+28="
$t.Cell(9, 3).Range.Text = "6+27="
$t.Cell(9, 4).Range.Text = "76-19="
$t.Cell(9, 5).Range.Text = "72-28="

$t.Cell(10, 1).Range.Text = "2+69="
$t.Cell(10, 2).Range.Text = "22+59="
$t.Cell(10, 3).Range.Text = "6+89="
$t.Cell(10, 4).Range.Text = "97-68="
$t.Cell(10, 5).Range.Text = "16+69="

$t.Cell(11, 1).Range.Text = "41-17="
$t.Cell(11, 2).Range.Text = "51-42="
$t.Cell(11, 3).Range.Text = "70-34="
$t.Cell(11, 4).Range.Text = "70-13="
$t.Cell(11, 5).Range.Text = "59+12="

$t.Cell(12, 1).Range.Text = "93-44="
$t.Cell(12, 2).Range.Text = "6+25="
$t.Cell(12, 3).Range.Text = "51-8="
$t.Cell(12, 4).Range.Text = "26+25="
$t.Cell(12, 5).Range.Text = "29+56="

$t.Cell(13, 1).Range.Text = "84-46="
$t.Cell(13, 2).Range.Text = "67+17="
$t.Cell(13, 3).Range.Text = "53+39="
$t.Cell(13, 4).Range.Text = "16+7="
$t.Cell(13, 5).Range.Text = "36+16="

$t.Cell(14, 1).Range.Text = "34-29="
$t.Cell(14, 2).Range.Text = "19+33="
$t.Cell(14, 3).Range.Text = "8+59="
$t.Cell(14, 4).Range.Text = "69+2="
$t.Cell(14, 5).Range.Text = "57+37="

$t.Cell(15, 1).Range.Text = "5+58="
$t.Cell(15, 2).Range.Text = "19+37="
$t.Cell(15, 3).Range.Text = "84-77="
$t.Cell(15, 4).Range.Text = "81-64="
$t.Cell(15, 5).Range.Text = "92-85="

$t.Cell(16, 1).Range.Text = "26+59="
$t.Cell(16, 2).Range.Text = "58+13="
$t.Cell(16, 3).Range.Text = "66+17="
$t.Cell(16, 4).Range.Text = "7+77="
$t.Cell(16, 5).Range.Text = "9+49="

$t.Cell(17, 1).Range.Text = "70-5="
$t.Cell(17, 2).Range.Text = "12+39="
$t.Cell(17, 3).Range.Text = "84-48="
$t.Cell(17, 4).Range.Text = "57-39="
$t.Cell(17, 5).Range.Text = "70-49="

$t.Cell(18, 1).Range.Text = "84-35="
$t.Cell(18, 2).Range.Text = "19+76="
$t.Cell(18, 3).Range.Text = "66+19="
$t.Cell(18, 4).Range.Text = "65+7="
$t.Cell(18, 5).Range.Text = "24+59="

$t.Cell(19, 1).Range.Text = "26-19="
$t.Cell(19, 2).Range.Text = "72-39="
$t.Cell(19, 3).Range.Text = "37+17="
$t.Cell(19, 4).Range.Text = "73-5="
$t.Cell(19, 5).Range.Text = "72-53="

$t.Cell(20, 1).Range.Text = "22+19="
$t.Cell(20, 2).Range.Text = "30-19="
$t.Cell(20, 3).Range.Text = "7+74="
$t.Cell(20, 4).Range.Text = "59+26="
$t.Cell(20, 5).Range.Text = "7+5="

